$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh Price (D) and Volume(1h) (E) columns with latest crypto data
$ws.Range("D2").Value = "28.068.85"
$ws.Range("E2").Value = "  -1.99%  "
$ws.Range("D3").Value = "1.799.96"
$ws.Range("E3").Value = "  +0.09%  "
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.32%  "
$ws.Range("D5").Value = "'316.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.89%  "
$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("D7").Value = "'0.5445"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.72%  "
$ws.Range("D8").Value = "'0.3789"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.19%  "
$ws.Range("D9").Value = "'0.07453"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.07%  "
$ws.Range("D10").Value = "'42.03"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.47%  "
$ws.Range("D11").Value = "'1.093"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.28%  "
$ws.Range("D12").Value = "'1.003"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.32%  "
$ws.Range("D13").Value = "'6.198"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("D14").Value = "'20.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.96%  "
$ws.Range("D15").Value = "'7.349"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.79%  "
$ws.Range("D16").Value = "1.800.69"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "'89.70"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.72%  "
$ws.Range("D18").Value = "'0.00001065"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").Value = "'0.06529"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.93%  "
$ws.Range("D20").Value = "'17.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.00%  "
$ws.Range("D21").Value = "'1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").Value = "'5.928"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "28.117.24"
$ws.Range("E23").Value = "  -1.83%  "
$ws.Range("D24").Value = "'11.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").Value = "'2.102"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("D26").Value = "'155.59"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.38%  "
$ws.Range("D27").Value = "'20.46"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("D28").Value = "2.011.60"
$ws.Range("E28").Value = "  +0.40%  "
$ws.Range("D29").Value = "'2.326"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.09%  "
$ws.Range("D30").Value = "'122.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.43%  "
$ws.Range("D31").Value = "'0.1120"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.38%  "
$ws.Range("D32").Value = "'1.115"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.14%  "
$ws.Range("D33").Value = "'3.666"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.45%  "
$ws.Range("D34").Value = "'5.552"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.60%  "
$ws.Range("D35").Value = "'0.06989"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +7.34%  "
$ws.Range("D36").Value = "'0.2224"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.15%  "
$ws.Range("D37").Value = "'0.02308"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").Value = "'5.065"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("D39").Value = "'8.456"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.11%  "
$ws.Range("D40").Value = "'11.17"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.52%  "
$ws.Range("D41").Value = "'0.6164"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.99%  "
$ws.Range("D44").Value = "'13.39"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.78%  "
$ws.Range("D45").Value = "'3.680"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("D46").Value = "'0.5746"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.84%  "
$ws.Range("D47").Value = "'124.65"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.47%  "
$ws.Range("D48").Value = "'1.184"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.22%  "
$ws.Range("D49").Value = "'1.920"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.82%  "
$ws.Range("D50").Value = "'0.06823"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.50%  "
$ws.Range("D51").Value = "'71.84"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.42%  "

# Rows 42 and 43: coin ranking reordered (TrustWalletToken <-> WEMIXTOKEN) with updated values
$ws.Range("B42").Value = "WEMIXTOKEN"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").Value = "'1.421"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.14%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "'1.170"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.11%  "
